# Add "Parallel Development" as a new bullet/benefit right after "Late Binding"
# on the "Benefits of Dependency Injection" slide.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$shp = $s.Shapes.Item("Content Placeholder 2")

$tr = $shp.TextFrame.TextRange

# Locate the paragraph that holds "Late Binding" so the new bullet is
# inserted right after it (and before the blank spacer paragraph).
$lateBindingPara = $null
for ($i = 1; $i -le $tr.Paragraphs().Count; $i++) {
    $para = $tr.Paragraphs($i, 1)
    if ($para.Text.Trim() -eq "Late Binding") {
        $lateBindingPara = $para
        break
    }
}

# Insert a new paragraph ("Parallel Development") right after "Late Binding".
# InsertAfter() on a paragraph-scoped range inherits that paragraph's run
# formatting (Arial, sz=2000) for the newly typed text; the leading `r`
# (carriage return) starts a new paragraph instead of appending inline.
$newPara = $lateBindingPara.InsertAfter("`rParallel Development")
